# Weekly update: a new price record (week of 2022-03-14) is inserted at the
# top of the data block (row 73), pushing all subsequent rows (old 73..122)
# down by one (new 74..123). The former last row (122, date 44160) becomes
# row 123 unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("73").Insert()

$ws.Range("A73").Value = 8
$ws.Range("B73").Value = "Terminal La Palmera de La Serena"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44634
$ws.Range("E73").Value = 4
$ws.Range("F73").Value = 100112001
$ws.Range("G73").Value = "Berenjena"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 480
$ws.Range("K73").Value = 8500
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = 8750
$ws.Range("N73").Value = "$/caja 50 unidades"
$ws.Range("O73").Value = "Región de Arica y Parinacota"
$ws.Range("P73").Value = 175
$ws.Range("Q73").Value = 50
$ws.Range("R73").Value = "Hortaliza"
